$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 holds the single data record in this test-data workbook.
# AV2 = FuncLoc, AX2 = Previous Doc -- update both to the new source values.
$ws.Range("AV2").Value = "ABCD330972"
$ws.Range("AX2").Value = "3050730281"
